# Generate Report for Handoff
# Adds a new handoff row (e76d9a39-679c-4e7d-9044-202a465d1302) to the
# "Overview" sheet and to each locale detail sheet ("zh-cn", "de-de"),
# mirroring the existing 3aafd78c... row already present in each sheet.

$wb = $excel.ActiveWorkbook

$newFileBase = "e76d9a39-679c-4e7d-9044-202a465d1302"
$newFileMd   = "$newFileBase.md"

# ---------------------------------------------------------------------
# Sheet "Overview": new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/425ef4b5ae0d5954614ea8acb24d3d27b481d81b/e2e/$newFileMd",
    "",
    "",
    $newFileMd
)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-37-13 14:37:09"

# ---------------------------------------------------------------------
# Sheet "zh-cn": new row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhHash = "9c164c373f71cf87186faa37af32418592217b0c"
$zhXlf  = "$newFileBase.$zhHash.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/425ef4b5ae0d5954614ea8acb24d3d27b481d81b/e2e/$newFileMd",
    "",
    "",
    $newFileMd
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/425ef4b5ae0d5954614ea8acb24d3d27b481d81b/e2e/$newFileMd",
    "",
    "",
    ".md"
)
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3212d225e8d7a9f3328a85b98493441103935ea1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf",
    "",
    "",
    $zhXlf
)
$wsZhCn.Range("E3").Value = "2016-03-13 14:37:06"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": new row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deHash = "9c164c373f71cf87186faa37af32418592217b0c"
$deXlf  = "$newFileBase.$deHash.de-de.xlf"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/425ef4b5ae0d5954614ea8acb24d3d27b481d81b/e2e/$newFileMd",
    "",
    "",
    $newFileMd
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/425ef4b5ae0d5954614ea8acb24d3d27b481d81b/e2e/$newFileMd",
    "",
    "",
    ".md"
)
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9d536afd391bef6cc1ad294551b38cf1ec821d7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf",
    "",
    "",
    $deXlf
)
$wsDeDe.Range("E3").Value = "2016-03-13 14:37:09"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"
